# The commit swaps the data of row 4 and row 5 in the "Artfynd" sheet,
# but only for columns A, B, E:J and Q:R. The remaining columns
# (C, D, P, S, T, U, V, W, Y, Z, AA, AB, AD, AE, AG, AT, AW, AX, AY)
# already hold identical values in both rows, so they do not need to move.
#
# Range.Copy(Destination) is used (instead of assigning .Value) so that
# Excel preserves the original cell types: numbers stay numbers, and the
# text values that look like numbers (e.g. "5"/"1" in column I) stay text,
# without introducing any new cell styles.
#
# Row 100 (well outside the worksheet's used range, A1:AY6) is used as
# scratch space to hold row 4's original values while row 5's values are
# copied into row 4.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Save row 4's changing cells into the scratch row (row 100).
$ws.Range("A4:B4").Copy($ws.Range("A100"))
$ws.Range("E4:J4").Copy($ws.Range("E100"))
$ws.Range("Q4:R4").Copy($ws.Range("Q100"))

# 2) Move row 5's changing cells into row 4.
$ws.Range("A5:B5").Copy($ws.Range("A4"))
$ws.Range("E5:J5").Copy($ws.Range("E4"))
$ws.Range("Q5:R5").Copy($ws.Range("Q4"))

# 3) Move the original row 4 values (held in the scratch row) into row 5.
$ws.Range("A100:B100").Copy($ws.Range("A5"))
$ws.Range("E100:J100").Copy($ws.Range("E5"))
$ws.Range("Q100:R100").Copy($ws.Range("Q5"))

# 4) Clean up the scratch row.
$ws.Range("A100:R100").Clear()
